# edit.ps1 - Applies the "Dark Matter" -> "Biology" content rewrite described
# by the provided diff, using Word COM-interop (Find/Replace) calls.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------------
Replace-Text "Exploring the Enigmatic Realm of Dark Matter" "Exploring the building blocks of life: Biology for the Curious"

# ---------------------------------------------------------------------------
# 2. Author name paragraph: "Dr. Amelia Vargas" (3 runs) -> "Madison Evans" (1 run)
# ---------------------------------------------------------------------------
Replace-Text "Dr. Amelia Vargas" "Madison Evans"

# ---------------------------------------------------------------------------
# 3. Email paragraph: "amelia.vargas@astrotech.edu" -> "madisonevans@hummingbirdmail.com"
# ---------------------------------------------------------------------------
Replace-Text "amelia" "madisonevans@hummingbirdmail"
Replace-Text "vargas@astrotech.edu" "com"

# ---------------------------------------------------------------------------
# 4. Body paragraph (first big block about dark matter -> biology)
# ---------------------------------------------------------------------------
Replace-Text "Within the vast cosmic canvas, there lies an enigmatic entity known as dark matter, an elusive substance that permeates the universe, yet remains shrouded in mystery" "Biology is the science that seeks to understand the intricate world of living organisms, encompassing everything from microscopic cells to majestic whales"

Replace-Text " Its existence has been inferred through its gravitational influence on visible matter, but its true nature continues to confound scientists" " It is a field that is both fascinating and challenging, with countless discoveries yet to be made"

Replace-Text " In this exploration, we delve into the enigmatic realm of dark matter, examining its intriguing properties, the methods employed to detect its presence, and the profound implications it holds for our understanding of the universe" " Delving into the realm of biology allows us to unravel the complexities of life, appreciate the interdependence of all living things, and gain a deeper understanding of our place in the universe"

Replace-Text "The existence of dark matter was first hinted at in the 1930s when observations of galaxy clusters revealed discrepancies between their expected and observed masses" "As we embark on this journey of exploration, we will peer into the inner workings of cells, the fundamental units of life"

Replace-Text " Since then, a multitude of astronomical observations have provided compelling evidence for its existence" " We will unravel the secrets of DNA, the blueprint for all living organisms, and witness the remarkable process of cellular respiration, which fuels the activities of life"

# Sentences 6 + "." + 7 collapse into a single run in the target.
Replace-Text " From the dynamics of galaxies and galaxy clusters to the gravitational lensing of light, the presence of dark matter is undeniable. Despite its profound influence, dark matter remains elusive to direct detection, leading to intense research efforts to unravel its secrets" " We will delve into the diverse ecosystems that blanket our planet, from lush rainforests to barren deserts, unraveling the delicate balance that sustains life in all its forms"

Replace-Text "Scientists employ various techniques to probe the enigmatic nature of dark matter" "Furthermore, we will investigate the intriguing world of genetics, exploring how genes influence traits and how variations in these genes drive evolution"

Replace-Text " One approach involves studying the gravitational effects it exerts on visible matter" " We will examine the fascinating relationship between organisms and their environment, deciphering how adaptations allow species to thrive in specific habitats"

# Sentences 10 + "." + 11 collapse into a single run in the target.
Replace-Text " By analyzing the motions of stars within galaxies and the dynamics of galaxy clusters, astronomers can infer the presence and distribution of dark matter. Additionally, sensitive detectors are employed in underground laboratories and space-based experiments to directly detect dark matter particles" " By unraveling these intricate connections, we gain insights into the remarkable resilience and adaptability of life on Earth"

# ---------------------------------------------------------------------------
# 5. Summary body paragraph
# ---------------------------------------------------------------------------
Replace-Text "Dark matter, an enigmatic entity permeating the universe, continues to captivate scientists with its elusive nature" "Biology is a captivating and multifaceted field that unveils the intricacies of life on Earth"

Replace-Text " Through gravitational observations, scientists have inferred its existence, but its true identity remains concealed" " It explores the fundamental building blocks of living organisms, the processes that drive their survival, and the interconnectedness of all living things"

# Sentences 3 + "." + 4 collapse into a single run in the target.
Replace-Text " The exploration of dark matter poses profound implications for our understanding of the universe, challenging conventional theories and prompting innovative research endeavors. As we delve deeper into the realm of dark matter, we unveil the mysteries that lie at the heart of our cosmic existence" " Through the study of biology, we gain a deeper understanding of the beauty and complexity of life and our role as stewards of this precious planet"

# ---------------------------------------------------------------------------
# 6. Add a trailing empty paragraph after the Summary body paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
